# Actualización automática 2025-06-06 08:00:08
#
# Records the new "SAL SOLUBLE" sale (2116.73) for client "ROCA REYNA PAUL
# DAVID" (advisor "LOZANO MOLINA TITO") and ripples the value through the
# three tracking sheets: the per-group sheet, the monthly sheet, and the
# monthly-compliance summary sheet (which also folds in two other
# already-pending sale amounts for FREGADEROS DE COCINA and PORCELANATO).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" -------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("N21").Value = 2116.73
$wsGrupo.Range("N29").Value = "1 de 27"

# --- Sheet 2: "VENTA MENSUAL" ----------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F21").Value = 2116.73
$wsMensual.Range("F29").Value = 2116.73
# column F ("junio") widens from 11 to 13 characters
# (COM ColumnWidth is offset by 5/6 from the stored OOXML <col width>)
$wsMensual.Columns.Item(6).ColumnWidth = 13 - 5/6

# --- Sheet 3: "CUMPLIMIENTO MENSUAL" ---------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# FREGADEROS DE COCINA (row 4)
$wsCumpl.Range("D4").Value = 560.03
$wsCumpl.Range("E4").Value = -309.398174579099
$wsCumpl.Range("F4").Value = 2.234472813097491

# PORCELANATO (row 16)
$wsCumpl.Range("D16").Value = 2695.68
$wsCumpl.Range("E16").Value = 10365.9
$wsCumpl.Range("F16").Value = 0.2063823825295255

# SAL SOLUBLE (row 18)
$wsCumpl.Range("D18").Value = 2116.73
$wsCumpl.Range("E18").Value = -916.73
$wsCumpl.Range("F18").Value = 1.763941666666667

# TOTAL (row 19)
$wsCumpl.Range("D19").Value = 5372.440000000001
$wsCumpl.Range("E19").Value = 18127.56093005039
$wsCumpl.Range("F19").Value = 0.2286144590373206

# columns D, E, F widen from 11/22/18 to 13/23/24 characters
$wsCumpl.Columns.Item(4).ColumnWidth = 13 - 5/6
$wsCumpl.Columns.Item(5).ColumnWidth = 23 - 5/6
$wsCumpl.Columns.Item(6).ColumnWidth = 24 - 5/6
